$d = $word.ActiveDocument

function Get-ParaIndexByText($substr) {
    $idx = 0
    foreach ($p in $d.Paragraphs) {
        $idx = $idx + 1
        if ($p.Range.Text -like $substr) {
            return $idx
        }
    }
    return -1
}

# Locate the paragraph holding the "{#table}" marker - it anchors the whole
# block that needs reworking.
$tableStartIdx = Get-ParaIndexByText("*{#table}*")

# 1) Drop the "{#table}" run entirely, leaving the (now run-less) paragraph
#    mark in place. Trim the trailing paragraph mark off the range first so
#    the mark itself is not deleted/merged with the next paragraph.
$tableMarkerRange = $d.Paragraphs.Item($tableStartIdx).Range
$tableMarkerRange = $d.Range($tableMarkerRange.Start, $tableMarkerRange.End - 1)
$tableMarkerRange.Delete()

# 2) The old markdown table-header paragraph ("| Q No | ... |") becomes the
#    new bold "PART- B " heading.
$headerIdx = $tableStartIdx + 1
$headerPara = $d.Paragraphs.Item($headerIdx)
$headerPara.Range.Text = "PART- B "
$headerPara.Range.Font.Bold = 1

# 3) Remove the now-redundant paragraphs that used to hold the markdown
#    table template (separator row, data row, "{/table}", blank spacer) plus
#    the original "PART- B" paragraph that used to follow them.
$deleteFromIdx = $headerIdx + 1
$deleteToIdx = $tableStartIdx + 6
$deleteRange = $d.Range($d.Paragraphs.Item($deleteFromIdx).Range.Start, $d.Paragraphs.Item($deleteToIdx).Range.End)
$deleteRange.Delete()

# 4) Remove the stray "//table 2" marker paragraph further down.
$table2Idx = Get-ParaIndexByText("*//table 2*")
if ($table2Idx -ne -1) {
    $d.Paragraphs.Item($table2Idx).Range.Delete()
}
